$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B, shifting the existing columns (B:E -> C:F)
$ws.Range("B1").EntireColumn.Insert()

# Set header values for the new layout: B..F
$ws.Range("B1").Value = "env"
$ws.Range("C1").Value = "id"
$ws.Range("D1").Value = "capacity"
$ws.Range("E1").Value = "material_nature"
$ws.Range("F1").Value = "kwargs"

# The insert left B1 without the header formatting (bold/border/center) that
# the other header cells carry, and F1 (brand-new cell) has none either.
# Restore/apply the shared header style by copying formats from neighboring
# header cells that already have it.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
